# Apply cryptocurrency price/volume updates per the diff.
# Values that look like plain numbers (e.g. "1.91") are written with a
# leading quote-prefix (') so Excel stores them as text, matching the
# original inlineStr cell type instead of silently converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "68.276.51"
$ws.Cells.Item(2, 5).Value = "  +1.47%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.737.55"
$ws.Cells.Item(3, 5).Value = "  +0.80%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.21%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'591.94"
$ws.Cells.Item(5, 5).Value = "  +0.84%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'167.18"
$ws.Cells.Item(6, 5).Value = "  +2.07%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "3.736.36"
$ws.Cells.Item(7, 5).Value = "  +0.74%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.06%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.519"
$ws.Cells.Item(9, 5).Value = "  +1.05%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +2.33%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +1.22%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.449"
$ws.Cells.Item(12, 5).Value = "  +1.19%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'0.0000260"
$ws.Cells.Item(13, 5).Value = "  +1.09%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'36.24"
$ws.Cells.Item(14, 5).Value = "  +2.52%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "4.361.53"
$ws.Cells.Item(15, 5).Value = "  +0.66%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "3.744.78"
$ws.Cells.Item(16, 5).Value = "  +0.98%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "68.226.61"
$ws.Cells.Item(17, 5).Value = "  +1.32%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'17.89"
$ws.Cells.Item(18, 5).Value = "  -1.47%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'7.01"
$ws.Cells.Item(19, 5).Value = "  +1.08%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +0.94%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'10.67"
$ws.Cells.Item(21, 5).Value = "  +1.73%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'466.04"
$ws.Cells.Item(22, 5).Value = "  +1.19%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'0.697"
$ws.Cells.Item(23, 5).Value = "  +0.96%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'0.0000149"
$ws.Cells.Item(24, 5).Value = "  +11.54%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'83.85"
$ws.Cells.Item(25, 5).Value = "  +2.03%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'2.18"
$ws.Cells.Item(26, 5).Value = "  +3.03%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'11.89"
$ws.Cells.Item(27, 5).Value = "  +0.91%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'10.06"
$ws.Cells.Item(28, 5).Value = "  +0.53%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.01%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "3.879.83"
$ws.Cells.Item(30, 5).Value = "  +0.59%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'2.76"
$ws.Cells.Item(31, 5).Value = "  -1.47%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'7.29"
$ws.Cells.Item(32, 5).Value = "  +0.93%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'29.83"
$ws.Cells.Item(33, 5).Value = "  +1.97%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.01%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'9.20"
$ws.Cells.Item(35, 5).Value = "  +3.71%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "3.689.36"
$ws.Cells.Item(37, 5).Value = "  +0.83%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.101"
$ws.Cells.Item(38, 5).Value = "  +0.68%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'3.44"
$ws.Cells.Item(39, 5).Value = "  +2.84%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +2.89%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.99%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'5.78"
$ws.Cells.Item(42, 5).Value = "  +2.38%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.15%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'43.89"
$ws.Cells.Item(45, 5).Value = "  +17.37%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.300"
$ws.Cells.Item(46, 5).Value = "  +0.06%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Stacks"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(47, 4).Value = "'1.91"
$ws.Cells.Item(47, 5).Value = "  +1.74%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "OKB"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(48, 4).Value = "'46.56"
$ws.Cells.Item(48, 5).Value = "  +3.71%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'8.45"
$ws.Cells.Item(49, 5).Value = "  +0.55%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Bittensor"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(50, 4).Value = "'389.54"
$ws.Cells.Item(50, 5).Value = "  +0.71%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Monero"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(51, 4).Value = "'144.26"
$ws.Cells.Item(51, 5).Value = "  +0.37%  "
